$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'265.26"
$ws.Range("D3").Value = "'22.68"
$ws.Range("D4").Value = "'6.231"
$ws.Range("D5").Value = "'0.06152"
$ws.Range("D6").Value = "'3.567"
$ws.Range("D7").Value = "'6.716"
$ws.Range("D8").Value = "'1.349"
$ws.Range("D9").Value = "'0.8271"
$ws.Range("D10").Value = "'0.01355"
$ws.Range("D11").Value = "'0.1594"
$ws.Range("D12").Value = "'0.08179"
$ws.Range("D13").Value = "'0.03398"
$ws.Range("D14").Value = "'0.03160"
$ws.Range("D15").Value = "'0.09232"
$ws.Range("D16").Value = "'3.902"
$ws.Range("D17").Value = "'0.001684"
$ws.Range("D18").Value = "'0.04786"
$ws.Range("D19").Value = "'0.006255"
$ws.Range("D20").Value = "'0.006308"
$ws.Range("D21").Value = "'0.001100"
$ws.Range("D22").Value = "'0.0001501"
$ws.Range("D23").Value = "'3.745"
$ws.Range("D25").Value = "'0.3352"
$ws.Range("D40").Value = "'0.04613"
$ws.Range("D41").Value = "'0.006979"
$ws.Range("D42").Value = "'0.1136"
$ws.Range("D43").Value = "'0.003239"
$ws.Range("D44").Value = "'0.01080"
$ws.Range("D45").Value = "'0.00006162"
$ws.Range("D47").Value = "'0.7709"
$ws.Range("E47").Value = "46CoinbaseStockTokenCOINWorstin24h"
$ws.Range("D48").Value = "'0.2059"
$ws.Range("D49").Value = "'0.00002103"
$ws.Range("E49").Value = "48CryptobidCoinCBC"
$ws.Range("D50").Value = "'0.01241"
